$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cartesian Velocity - Linear")

$ws.Range("B3").Value = 0.3136197481143909
$ws.Range("C3").Value = 0.3516260814693335
$ws.Range("D3").Value = 0.325169621440986
$ws.Range("E3").Value = 0.2641859620991784

$ws.Range("B4").Value = 0.9921422580253202
$ws.Range("C4").Value = 0.999768409175526
$ws.Range("D4").Value = 0.9951598255184436
$ws.Range("E4").Value = 0.9918503643734413

$ws.Range("B5").Value = 0.4416510353234783
$ws.Range("C5").Value = 0.529979252502416
$ws.Range("D5").Value = 0.6418049323436407
$ws.Range("E5").Value = 0.8334310143006941
